# Change the "From" value of rule R40 (row 10) on the Rules sheet
# from 18 to 1, as captured in the commit's restored revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
